$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# STEP 1: Insert a new column before column B (shifts old B..L -> C..M across the whole sheet).
# This naturally fixes row 2 (title merge A2:M2), row 10 (new column for SUCURSAL),
# and gives rows 5-8 their extra trailing column.
$ws.Columns("B").Insert()

# Give the freshly inserted column B the same width the old column B used to have
# (closest achievable value to the original 17.6640625 stored width).
$ws.Columns("B").ColumnWidth = 16.75

# STEP 2: Row 10 header - put "SUCURSAL" into the newly inserted B10 cell.
$ws.Range("B10").Value = "SUCURSAL"

# STEP 3: Rows 5-8 - the value-merge (originally B:D, now shifted to C:E) needs to be widened
# to include column B as well (B:E), matching the new field width.
foreach ($r in 5..8) {
    $ws.Range("C" + $r + ":E" + $r).UnMerge()
    $ws.Range("D" + $r).Copy($ws.Range("B" + $r))
    $ws.Range("B" + $r + ":E" + $r).Merge()
}

# STEP 4: Row 4 - restructure into two label/value pairs:
#   A4 = "EMPRESA:"            B4:E4 = value field
#   F4 = "ESTABLECIMEINTO :"   G4:I4 = value field
$ws.Range("C4:E4").UnMerge()

# Copy the old A4 (label style, currently still "ESTABLECIMEINTO :") into F4 before we overwrite A4.
$ws.Range("A4").Copy($ws.Range("F4"))

# Now change A4's text to the new label.
$ws.Range("A4").Value = "EMPRESA:"

# Fix up the B4 style (currently inherited bold/label style from the insert) to the plain value style.
$ws.Range("D4").Copy($ws.Range("B4"))

# Build the plain value-style cells for the new ESTABLECIMEINTO value field (G4:I4).
$ws.Range("D4").Copy($ws.Range("G4"))
$ws.Range("D4").Copy($ws.Range("H4"))
$ws.Range("D4").Copy($ws.Range("I4"))

# Merge the two value fields on row 4.
$ws.Range("B4:E4").Merge()
$ws.Range("G4:I4").Merge()

# Restore the originally selected cell.
$ws.Range("B10").Select()

Write-Host "done"
